$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.987.76"
$ws.Range("D3").Value = "2.477.01"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.75"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.18"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").Value = "2.476.58"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.84"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "68.897.41"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.67"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.64"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.00"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.90"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.97"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.28"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("D26").Value = "0.0₃0827"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "429.37"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.15"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.02"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.93"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.44"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.89"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.486"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.568"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.99"
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.96"
$ws.Range("E51").Value = "  -2.24%  "
